$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.675.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.08%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.852.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.01%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.36%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4259"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.02%  "

$ws.Range("E8").Value = "  +0.87%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.58"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07301"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.36%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8755"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.64"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.12%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.828.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.95%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.520"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.55%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.318"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06908"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.50%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.42%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "79.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.50%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009016"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9984"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.683.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.94%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.976"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.061.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.59%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.981"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.72%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "121.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +10.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.268"
$ws.Range("D30").Style = "Normal"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.863"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08909"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.69%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7697"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.07%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.528"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.93%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.964"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.77%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.109"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05385"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.94%  "

$ws.Range("E38").Value = "  +0.83%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01940"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.88%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.822"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.78%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5060"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.70%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.814"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.38%  "

$ws.Range("E43").Value = "  +1.54%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.370"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.62%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06540"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "104.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4668"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.60%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9987"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.14%  "

$ws.Range("E50").Value = "  -0.14%  "

$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.756"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.00%  "
